$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate current row 8 (Fecha 2021-07-20, Volumen 60, etc.) down into the
# new row 9, carrying formatting (incl. the date number format on column D)
# along with it.
$ws.Range("A8:T8").Copy($ws.Range("A9:T9"))

# Now update row 8 in place with the new weekly figures.
$ws.Range("D8").Value = 44491
$ws.Range("M8").Value = 180
$ws.Range("N8").Value = 9000
$ws.Range("O8").Value = 9000
$ws.Range("P8").Value = 9000
$ws.Range("Q8").Value = '$/caja 14 kilos empedrada'
$ws.Range("S8").Value = 643
